$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the underlying data values that drive both charts.
$ws.Range("C2").Value = 0.27229999999999999
$ws.Range("D2").Value = 0.28270000000000001
$ws.Range("L2").Value = 45.131
$ws.Range("M2").Value = 41.958799999999997

# Restore the active sheet selection to D3 as recorded after the edit.
$ws.Activate()
$ws.Range("D3").Select()
